$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Budget feature: the expense table grew from A1:C6 to A1:C11. Several new
# expense rows were inserted after the header, the "food"/"Travel" rows were
# replaced with fresh entries, while "cat" and "Tuition fee" were kept (now
# further down the table) and the last "Food" row stayed the same.

# First, remember the date-format style used by the existing date column
# (C2 currently carries it) so every new/ shifted row in column C keeps the
# same formatting.
$ws.Cells.Item(2, 3).Copy() | Out-Null

$data = @(
    @("fooD", 2000, 45908.250231481485),
    @("Travel", 100, 45907.250231481485),
    @("Food", 2000, 45906.250231481485),
    @("Food", 1200, 45905.250231481485),
    @("Games", 1000, 45904.250231481485),
    @("Extra", 1000, 45903.250231481485),
    @("Food", 1500, 45903.250231481485),
    @("cat", 30000, 45883.250231481485),
    @("Tuition fee", 90000, 45873.250231481485),
    @("Food", 10000, 45870.250231481485)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 3).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$excel.CutCopyMode = $false
